# Data Tables & Report Types
# Duplicate the "Results3" report sheet to create a new "Results4" sheet
# (2 additional report categories: source and Meter use the same table
# layout/structure, so the new sheet is produced from the existing one).

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("Results3")
$lastSheet   = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy the existing report sheet to the end of the workbook and rename it.
$sourceSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Results4"

# Keep the originally active sheet selected (adding the new report sheet
# should not change which tab is active).
$wb.Worksheets.Item(1).Activate()
